# Generate Report for Handoff
# Inserts a new tracked file "bbc70313-0b99-49cb-bd14-5bb60855b5e6.md" (status
# "Ready for handoff") ahead of the existing "e713a193-...md" row on every
# sheet (Overview, zh-cn, de-de), pushing that row and the trailing
# ".localization-config" row down by one. Hyperlinks are rebuilt from scratch
# on every sheet so the ref attributes line up with the new row numbers.

$wb = $excel.ActiveWorkbook

$newFile = "bbc70313-0b99-49cb-bd14-5bb60855b5e6.md"
$newStatus = "Ready for handoff"
$newHash = "7612dee54b7921c13ec218183cb04a1d019c48de"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Duplicate row 6 (the "e713a193..." row) in place so the new row inherits
# identical cell formatting, then push everything from the old row 6 down.
$ws1.Rows.Item(6).Copy()
$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value2 = $newFile
$ws1.Range("B6").Value2 = $newStatus
$ws1.Range("C6").Value2 = $newStatus

# Rebuild every hyperlink on the sheet (row numbers below A6 shifted by +1).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md", "", "", "030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/5f81dd56-9354-4281-aefc-7062931e76a6.md", "", "", "5f81dd56-9354-4281-aefc-7062931e76a6.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/91151bbfe283826399f6e971c0b4643395bdfd2e/e2e/7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md", "", "", "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/cf9ac2ce-8423-4665-a7ba-d02d0df06863.md", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/bbc70313-0b99-49cb-bd14-5bb60855b5e6.md", "", "", $newFile) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/935aa4c034d16f274e4a7110fc617ab92a93172a/e2e/e713a193-f786-4add-a55d-bb609b3d2b0f.md", "", "", "e713a193-f786-4add-a55d-bb609b3d2b0f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(6).Copy()
$ws2.Rows.Item(6).Insert()

$ws2.Range("A6").Value2 = $newFile
$ws2.Range("B6").Value2 = $newStatus
$ws2.Range("C6").Value2 = "$newFile.$newHash.zh-cn.xlf"
$ws2.Range("D6").Value2 = "2016-03-09 05:32:25"
$ws2.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H6").Value2 = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md", "", "", "030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b894eb86246409b9e4a158dfa9d125b4c3c36625/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/030c1e0f-1a18-4e77-a619-0a025f4dbaa8.b03d696bb6f815287c9b70fe5ad9ea1b51208234.zh-cn.xlf", "", "", "030c1e0f-1a18-4e77-a619-0a025f4dbaa8.b03d696bb6f815287c9b70fe5ad9ea1b51208234.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/5f81dd56-9354-4281-aefc-7062931e76a6.md", "", "", "5f81dd56-9354-4281-aefc-7062931e76a6.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b894eb86246409b9e4a158dfa9d125b4c3c36625/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5f81dd56-9354-4281-aefc-7062931e76a6.e37b314ad4428b0f90f296d80ac74f75fb2b761f.zh-cn.xlf", "", "", "5f81dd56-9354-4281-aefc-7062931e76a6.e37b314ad4428b0f90f296d80ac74f75fb2b761f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/91151bbfe283826399f6e971c0b4643395bdfd2e/e2e/7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md", "", "", "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4e32cd213c3c9230c11d1c4b444f09a436a6ead/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.b8c3a17b9e32384997bc4977e2ed0ac174c847b9.zh-cn.xlf", "", "", "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.b8c3a17b9e32384997bc4977e2ed0ac174c847b9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/cf9ac2ce-8423-4665-a7ba-d02d0df06863.md", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3aa4566ee0e8bced4fd8ea1b1095dc7d50c5e048/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.zh-cn.xlf", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/94157614ada4a1658afb5f44e8601f0b73352f57/e2e/cf9ac2ce-8423-4665-a7ba-d02d0df06863.md", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d423b9f17d9c3f8e30ba58bdf28cf62c50f69f8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.zh-cn.xlf", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/bbc70313-0b99-49cb-bd14-5bb60855b5e6.md", "", "", $newFile) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newFile.$newHash.zh-cn.xlf", "", "", "$newFile.$newHash.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/935aa4c034d16f274e4a7110fc617ab92a93172a/e2e/e713a193-f786-4add-a55d-bb609b3d2b0f.md", "", "", "e713a193-f786-4add-a55d-bb609b3d2b0f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf3253524abc984fd333fcd76e1e0e8fa4f1a9e3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e713a193-f786-4add-a55d-bb609b3d2b0f.7841ff6c759e5ae0f2cabd004e50ac176a47cd3b.zh-cn.xlf", "", "", "e713a193-f786-4add-a55d-bb609b3d2b0f.7841ff6c759e5ae0f2cabd004e50ac176a47cd3b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(6).Copy()
$ws3.Rows.Item(6).Insert()

$ws3.Range("A6").Value2 = $newFile
$ws3.Range("B6").Value2 = $newStatus
$ws3.Range("C6").Value2 = "$newFile.$newHash.de-de.xlf"
$ws3.Range("D6").Value2 = "2016-03-09 05:32:33"
$ws3.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H6").Value2 = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md", "", "", "030c1e0f-1a18-4e77-a619-0a025f4dbaa8.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fb024b046bf8684bcbdf7974f9225437e70950b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/030c1e0f-1a18-4e77-a619-0a025f4dbaa8.b03d696bb6f815287c9b70fe5ad9ea1b51208234.de-de.xlf", "", "", "030c1e0f-1a18-4e77-a619-0a025f4dbaa8.b03d696bb6f815287c9b70fe5ad9ea1b51208234.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5253f095428b8ebaa769b1668531039ec14b7bd/e2e/5f81dd56-9354-4281-aefc-7062931e76a6.md", "", "", "5f81dd56-9354-4281-aefc-7062931e76a6.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fb024b046bf8684bcbdf7974f9225437e70950b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5f81dd56-9354-4281-aefc-7062931e76a6.e37b314ad4428b0f90f296d80ac74f75fb2b761f.de-de.xlf", "", "", "5f81dd56-9354-4281-aefc-7062931e76a6.e37b314ad4428b0f90f296d80ac74f75fb2b761f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/91151bbfe283826399f6e971c0b4643395bdfd2e/e2e/7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md", "", "", "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bfc8d598aad35c345f78a59966fb8ca5e1526cbd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.b8c3a17b9e32384997bc4977e2ed0ac174c847b9.de-de.xlf", "", "", "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.b8c3a17b9e32384997bc4977e2ed0ac174c847b9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/cf9ac2ce-8423-4665-a7ba-d02d0df06863.md", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a86bb46b68f6d7f7bf91f5ad756360c0648dbb94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.de-de.xlf", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f942b3db5fcdf40486205413d7c33319b3a0fa03/e2e/cf9ac2ce-8423-4665-a7ba-d02d0df06863.md", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6faf5be841a6cfda82b64345e225213249ac038c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.de-de.xlf", "", "", "cf9ac2ce-8423-4665-a7ba-d02d0df06863.c05b13a9a0d2aa4987a6007f3b12ae163a13e739.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/e2e/bbc70313-0b99-49cb-bd14-5bb60855b5e6.md", "", "", $newFile) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newFile.$newHash.de-de.xlf", "", "", "$newFile.$newHash.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/935aa4c034d16f274e4a7110fc617ab92a93172a/e2e/e713a193-f786-4add-a55d-bb609b3d2b0f.md", "", "", "e713a193-f786-4add-a55d-bb609b3d2b0f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83be85b1c05b67310707a274e06371b623abc393/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e713a193-f786-4add-a55d-bb609b3d2b0f.7841ff6c759e5ae0f2cabd004e50ac176a47cd3b.de-de.xlf", "", "", "e713a193-f786-4add-a55d-bb609b3d2b0f.7841ff6c759e5ae0f2cabd004e50ac176a47cd3b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/b22808a5d042605f7fd769d8f8fcfc0ec36dcf2e/.localization-config", "", "", ".localization-config") | Out-Null

$wb.Save()
